$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 47

# Column A: plain text date string (matches the other "MM/DD/YYYY" text rows,
# e.g. A10, A37..A46), not an Excel date serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "09/07/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.0004494500000000005
$ws.Cells.Item($row, 3).Value = 111247.079764156
$ws.Cells.Item($row, 4).Value = 50
